$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly record needs to be inserted as row 14 (pushing the existing
# rows 14-31 down to 15-32). The new row repeats the same constant values
# (market, region, category, etc.) as the surrounding rows, but carries its
# own date / price figures.

$ws.Rows.Item(14).Insert()

$ws.Cells.Item(14, 1).Value = 6
$ws.Cells.Item(14, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(14, 3).Value = "Metropolitana"
$ws.Cells.Item(14, 4).Value = 44775
$ws.Cells.Item(14, 4).NumberFormat = $ws.Cells.Item(15, 4).NumberFormat
$ws.Cells.Item(14, 5).Value = 13
$ws.Cells.Item(14, 6).Value = 100112035
$ws.Cells.Item(14, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(14, 8).Value = "Sin especificar"
$ws.Cells.Item(14, 9).Value = "Primera"
$ws.Cells.Item(14, 10).Value = 250
$ws.Cells.Item(14, 11).Value = 18000
$ws.Cells.Item(14, 12).Value = 20000
$ws.Cells.Item(14, 13).Value = 19200
$ws.Cells.Item(14, 14).Value = "$/malla 15 kilos"
$ws.Cells.Item(14, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(14, 16).Value = 1280
$ws.Cells.Item(14, 17).Value = 15
$ws.Cells.Item(14, 18).Value = "Hortaliza"
